$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet logs weekly wholesale-market Ciboulette prices; the commit adds a
# new week's data. Insert two new rows at row 442, pushing existing rows
# 442+ down to 444+ (dimension grows from R483 to R485).
$ws.Rows.Item(442).Insert()
$ws.Rows.Item(442).Insert()

# Fill new row 442 with data (Fecha 44769 = 2022-07-27, Calidad Primera)
$ws.Cells.Item(442, 1).Value = 6
$ws.Cells.Item(442, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(442, 3).Value = 'Metropolitana'
$ws.Cells.Item(442, 4).Value = 44769
$ws.Cells.Item(442, 5).Value = 13
$ws.Cells.Item(442, 6).Value = 100112039
$ws.Cells.Item(442, 7).Value = 'Ciboulette'
$ws.Cells.Item(442, 8).Value = 'Sin especificar'
$ws.Cells.Item(442, 9).Value = 'Primera'
$ws.Cells.Item(442, 10).Value = 290
$ws.Cells.Item(442, 11).Value = 2500
$ws.Cells.Item(442, 12).Value = 2500
$ws.Cells.Item(442, 13).Value = 2500
$ws.Cells.Item(442, 14).Value = '$/docena de atados'
$ws.Cells.Item(442, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(442, 16).Value = 833
$ws.Cells.Item(442, 17).Value = 3
$ws.Cells.Item(442, 18).Value = 'Hortaliza'

# Fill new row 443 with data (Fecha 44769 = 2022-07-27, Calidad Segunda)
$ws.Cells.Item(443, 1).Value = 6
$ws.Cells.Item(443, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(443, 3).Value = 'Metropolitana'
$ws.Cells.Item(443, 4).Value = 44769
$ws.Cells.Item(443, 5).Value = 13
$ws.Cells.Item(443, 6).Value = 100112039
$ws.Cells.Item(443, 7).Value = 'Ciboulette'
$ws.Cells.Item(443, 8).Value = 'Sin especificar'
$ws.Cells.Item(443, 9).Value = 'Segunda'
$ws.Cells.Item(443, 10).Value = 220
$ws.Cells.Item(443, 11).Value = 2000
$ws.Cells.Item(443, 12).Value = 2000
$ws.Cells.Item(443, 13).Value = 2000
$ws.Cells.Item(443, 14).Value = '$/docena de atados'
$ws.Cells.Item(443, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(443, 16).Value = 667
$ws.Cells.Item(443, 17).Value = 3
$ws.Cells.Item(443, 18).Value = 'Hortaliza'
